$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    "Match TV",
    "MBC (Middle East Broadcasting Center)",
    "MBC (Munhwa Broadcasting Corporation)",
    "MEASAT Broadcast Network Systems",
    "Media Broadcast",
    "Media Nusantara Citra",
    "MediaCorp TV",
    "Mediapro",
    "Mediaset",
    "Mediawan Thematics",
    "MEGA TV Greece",
    "Megacable (Mega movil)",
    "Megogo",
    "Middle East Broadcasting Center",
    "Millicom",
    "Mola",
    "Movistar",
    "MPC (Moving Pictures)",
    "MTN",
    "MTS (Mobile TeleSystems)",
    "MTS (Telekom Srbija)"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
